# Update the "Förändrad" (column C) date value for rows 2-16 from
# serial date 45185 (2023-09-16) to serial date 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$oldSerial = 45185
$newSerial = 45204

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
